# Adjust Investment Summary table column widths for better formatting
#
# The three comparison/summary tables (slide 2 "Why This Solution?",
# slide 3 "Business Value - Financial Impact", slide 4 "Risk Mitigation")
# each get their placeholder content cleared out and their last table
# column nudged by 1 EMU so the columns add back up to a slightly wider
# overall table extent (the frame's <a:ext> is re-derived automatically
# from the table geometry).

$p = $ppt.ActivePresentation

function Clear-Table($tbl) {
    $rowCount = $tbl.Rows.Count
    $colCount = $tbl.Columns.Count
    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $tbl.Cell($r, $c).Shape.TextFrame.TextRange.Text = ""
        }
    }
}

# --- Slide 2: "Traditional Approach" vs "Our Solution" table (2 cols) ---
$s2 = $p.Slides.Item(2)
$tbl2 = $s2.Shapes.Item(3).Table
Clear-Table $tbl2
$tbl2.Columns.Item(2).Width = 342.950157480315

# --- Slide 3: "Metric" / "Value" investment summary table (2 cols) ---
$s3 = $p.Slides.Item(3)
$tbl3 = $s3.Shapes.Item(3).Table
Clear-Table $tbl3
$tbl3.Columns.Item(2).Width = 342.950157480315

# --- Slide 4: "Risk" / "Mitigation Strategy" / "Success Probability" table (3 cols) ---
$s4 = $p.Slides.Item(4)
$tbl4 = $s4.Shapes.Item(3).Table
Clear-Table $tbl4
$tbl4.Columns.Item(3).Width = 228.633464566929
